$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed metric values for rows 1-35 (unchanged row positions)
$ws.Cells.Item(1,2).Value = 0.497
$ws.Cells.Item(1,3).Value = 0.002
$ws.Cells.Item(2,2).Value = 0.865
$ws.Cells.Item(2,3).Value = 0.003
$ws.Cells.Item(3,2).Value = 0.835
$ws.Cells.Item(3,3).Value = 0.003
$ws.Cells.Item(4,2).Value = 0.437
$ws.Cells.Item(5,2).Value = 0.442
$ws.Cells.Item(6,2).Value = 1.22
$ws.Cells.Item(6,3).Value = 0.004
$ws.Cells.Item(6,4).Value = 0.629
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(7,2).Value = 1.1
$ws.Cells.Item(7,3).Value = 0.004
$ws.Cells.Item(7,4).Value = 0.777
$ws.Cells.Item(7,5).Value = 0.857
$ws.Cells.Item(8,2).Value = 1.316
$ws.Cells.Item(8,4).Value = 0.97
$ws.Cells.Item(9,2).Value = 1.114
$ws.Cells.Item(9,3).Value = 0.004
$ws.Cells.Item(9,4).Value = 0.755
$ws.Cells.Item(10,2).Value = 1.253
$ws.Cells.Item(10,4).Value = 0.754
$ws.Cells.Item(10,5).Value = 0.857
$ws.Cells.Item(11,2).Value = 0.85
$ws.Cells.Item(11,3).Value = 0.003
$ws.Cells.Item(12,2).Value = 1.086
$ws.Cells.Item(12,3).Value = 0.004
$ws.Cells.Item(13,2).Value = 0.831
$ws.Cells.Item(13,3).Value = 0.003
$ws.Cells.Item(13,4).Value = 0.8100000000000001
$ws.Cells.Item(14,2).Value = 0.9350000000000001
$ws.Cells.Item(15,2).Value = 0.9330000000000001
$ws.Cells.Item(15,4).Value = 0.865
$ws.Cells.Item(16,2).Value = 0.713
$ws.Cells.Item(17,2).Value = 0.699
$ws.Cells.Item(18,2).Value = 1.01
$ws.Cells.Item(18,4).Value = 0.907
$ws.Cells.Item(19,2).Value = 0.655
$ws.Cells.Item(20,2).Value = 0.988
$ws.Cells.Item(20,4).Value = 0.827
$ws.Cells.Item(21,2).Value = 0.521
$ws.Cells.Item(21,3).Value = 0.002
$ws.Cells.Item(22,2).Value = 1.119
$ws.Cells.Item(22,3).Value = 0.004
$ws.Cells.Item(22,4).Value = 0.716
$ws.Cells.Item(23,2).Value = 1.126
$ws.Cells.Item(23,3).Value = 0.004
$ws.Cells.Item(23,4).Value = 0.864
$ws.Cells.Item(24,2).Value = 1.157
$ws.Cells.Item(24,3).Value = 0.004
$ws.Cells.Item(24,4).Value = 0.924
$ws.Cells.Item(25,2).Value = 1.167
$ws.Cells.Item(25,3).Value = 0.004
$ws.Cells.Item(25,4).Value = 0.931
$ws.Cells.Item(26,2).Value = 0.9409999999999999
$ws.Cells.Item(26,4).Value = 0.531
$ws.Cells.Item(26,5).Value = 0.714
$ws.Cells.Item(27,2).Value = 1.174
$ws.Cells.Item(27,3).Value = 0.004
$ws.Cells.Item(27,4).Value = 0.929
$ws.Cells.Item(28,2).Value = 1.185
$ws.Cells.Item(28,3).Value = 0.004
$ws.Cells.Item(28,4).Value = 0.923
$ws.Cells.Item(29,2).Value = 0.896
$ws.Cells.Item(29,3).Value = 0.003
$ws.Cells.Item(30,2).Value = 1.173
$ws.Cells.Item(30,3).Value = 0.004
$ws.Cells.Item(30,4).Value = 0.88
$ws.Cells.Item(31,2).Value = 0.658
$ws.Cells.Item(31,4).Value = 0.832
$ws.Cells.Item(32,2).Value = 0.73
$ws.Cells.Item(32,4).Value = 0.79
$ws.Cells.Item(33,2).Value = 1.034
$ws.Cells.Item(33,4).Value = 0.889
$ws.Cells.Item(33,5).Value = 0.571
$ws.Cells.Item(34,2).Value = 1.233
$ws.Cells.Item(34,3).Value = 0.004
$ws.Cells.Item(34,4).Value = 0.853
$ws.Cells.Item(34,5).Value = 0.714
$ws.Cells.Item(35,2).Value = 1.271
$ws.Cells.Item(35,3).Value = 0.004
$ws.Cells.Item(35,4).Value = 0.902

# Insert a new row at position 36, shifting old rows 36-44 down to 37-45
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 (FY_4.png)
$ws.Cells.Item(36,1).Value = 'FY_4.png'
$ws.Cells.Item(36,2).Value = 1.167
$ws.Cells.Item(36,3).Value = 0.004
$ws.Cells.Item(36,4).Value = 0.841
$ws.Cells.Item(36,5).Value = 0.714
$ws.Cells.Item(36,6).Value = 'Fanny Yusuf'
$ws.Cells.Item(36,7).Value = 'Benar'

# Update shifted rows 37-45 with their final values
# Row 37 (TO_1.png)
$ws.Cells.Item(37,1).Value = 'TO_1.png'
$ws.Cells.Item(37,2).Value = 0.822
$ws.Cells.Item(37,3).Value = 0.003
$ws.Cells.Item(37,4).Value = 0.806
$ws.Cells.Item(37,5).Value = 1
$ws.Cells.Item(37,6).Value = 'Tiara Oktavian'
$ws.Cells.Item(37,7).Value = 'Benar'
# Row 38 (TO_2.png)
$ws.Cells.Item(38,1).Value = 'TO_2.png'
$ws.Cells.Item(38,2).Value = 0.914
$ws.Cells.Item(38,3).Value = 0.003
$ws.Cells.Item(38,4).Value = 0.867
$ws.Cells.Item(38,5).Value = 1
$ws.Cells.Item(38,6).Value = 'Tiara Oktavian'
$ws.Cells.Item(38,7).Value = 'Benar'
# Row 39 (TO_3.png)
$ws.Cells.Item(39,1).Value = 'TO_3.png'
$ws.Cells.Item(39,2).Value = 0.847
$ws.Cells.Item(39,3).Value = 0.003
$ws.Cells.Item(39,4).Value = 0.858
$ws.Cells.Item(39,5).Value = 1
$ws.Cells.Item(39,6).Value = 'Tiara Oktavian'
$ws.Cells.Item(39,7).Value = 'Benar'
# Row 40 (TO_4.png)
$ws.Cells.Item(40,1).Value = 'TO_4.png'
$ws.Cells.Item(40,2).Value = 2.245
$ws.Cells.Item(40,3).Value = 0.007
$ws.Cells.Item(40,4).Value = 0.506
$ws.Cells.Item(40,5).Value = 1
$ws.Cells.Item(40,6).Value = 'Tiara Oktavian'
$ws.Cells.Item(40,7).Value = 'Benar'
# Row 41 (TO_5.png)
$ws.Cells.Item(41,1).Value = 'TO_5.png'
$ws.Cells.Item(41,2).Value = 1.981
$ws.Cells.Item(41,3).Value = 0.006
$ws.Cells.Item(41,4).Value = 0.463
$ws.Cells.Item(41,5).Value = 1
$ws.Cells.Item(41,6).Value = 'Tiara Oktavian'
$ws.Cells.Item(41,7).Value = 'Benar'
# Row 42 (TD_1.png)
$ws.Cells.Item(42,1).Value = 'TD_1.png'
$ws.Cells.Item(42,2).Value = 1.47
$ws.Cells.Item(42,3).Value = 0.005
$ws.Cells.Item(42,4).Value = 0.421
$ws.Cells.Item(42,5).Value = 0.286
$ws.Cells.Item(42,6).Value = 'Tidak Diketahui'
$ws.Cells.Item(42,7).Value = 'Benar'
# Row 43 (TD_2.png)
$ws.Cells.Item(43,1).Value = 'TD_2.png'
$ws.Cells.Item(43,2).Value = 1.507
$ws.Cells.Item(43,3).Value = 0.005
$ws.Cells.Item(43,4).Value = 0.372
$ws.Cells.Item(43,5).Value = 0.571
$ws.Cells.Item(43,6).Value = 'Muhammad Iqbal Baqi'
$ws.Cells.Item(43,7).Value = 'Salah'
# Row 44 (TD_3.png)
$ws.Cells.Item(44,1).Value = 'TD_3.png'
$ws.Cells.Item(44,2).Value = 1.079
$ws.Cells.Item(44,3).Value = 0.004
$ws.Cells.Item(44,4).Value = 0.738
$ws.Cells.Item(44,5).Value = 0.429
$ws.Cells.Item(44,6).Value = 'Tidak Diketahui'
$ws.Cells.Item(44,7).Value = 'Benar'
# Row 45 (TD_4.png)
$ws.Cells.Item(45,1).Value = 'TD_4.png'
$ws.Cells.Item(45,2).Value = 0.992
$ws.Cells.Item(45,3).Value = 0.003
$ws.Cells.Item(45,4).Value = 0.703
$ws.Cells.Item(45,5).Value = 0.286
$ws.Cells.Item(45,6).Value = 'Tidak Diketahui'
$ws.Cells.Item(45,7).Value = 'Benar'
